$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 532
$ws.Range("F6").Value = 954
$ws.Range("F7").Value = 190
$ws.Range("F9").Value = 1019
$ws.Range("F10").Value = 811
$ws.Range("F11").Value = 237
$ws.Range("F14").Value = 818
$ws.Range("F15").Value = 276
$ws.Range("G16").Value = "已售罄"
$ws.Range("F17").Value = 500
$ws.Range("F18").Value = 1326
$ws.Range("F21").Value = 1172
$ws.Range("F22").Value = 2856
$ws.Range("F23").Value = 1401
$ws.Range("F24").Value = 697
$ws.Range("F25").Value = 189
$ws.Range("F26").Value = 1269
$ws.Range("F28").Value = 1010
$ws.Range("F29").Value = 355
$ws.Range("F30").Value = 3053
$ws.Range("F31").Value = 591
$ws.Range("F32").Value = 533
$ws.Range("F33").Value = 1388

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 76
$ws.Range("F4").Value = 364
$ws.Range("F10").Value = 153
$ws.Range("F12").Value = 14

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 734

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 734
$ws.Range("F4").Value = 532
$ws.Range("F5").Value = 76
$ws.Range("F9").Value = 364
$ws.Range("F12").Value = 954
$ws.Range("F13").Value = 190
$ws.Range("F16").Value = 1019
$ws.Range("F17").Value = 811
$ws.Range("F18").Value = 237
$ws.Range("F22").Value = 153
$ws.Range("F25").Value = 14
$ws.Range("F26").Value = 818
$ws.Range("F27").Value = 276
$ws.Range("G28").Value = "已售罄"
$ws.Range("F29").Value = 500
$ws.Range("F30").Value = 1326
$ws.Range("F33").Value = 1172
$ws.Range("F34").Value = 2856
$ws.Range("F35").Value = 1401
$ws.Range("F36").Value = 697
$ws.Range("F37").Value = 189
$ws.Range("F38").Value = 1269
$ws.Range("F42").Value = 1010
$ws.Range("F43").Value = 355
$ws.Range("F44").Value = 3053
$ws.Range("F45").Value = 591
$ws.Range("F46").Value = 533
$ws.Range("F47").Value = 1388

